# Applies the "cryptos list" price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Order follows the sheet (rows 2..51, cols B..E).
$updates = [ordered]@{
    "D2" = "60.821.05"
    "E2" = "  +0.08%  "
    "D3" = "3.363.24"
    "E3" = "  -0.69%  "
    "E4" = "  +0.01%  "
    "D5" = "569.53"
    "E5" = "  -0.19%  "
    "D6" = "138.49"
    "E6" = "  -2.38%  "
    "E7" = "  -0.01%  "
    "E8" = "  -0.66%  "
    "D9" = "7.61"
    "E9" = "  +1.30%  "
    "E10" = "  -2.51%  "
    "D11" = "0.379"
    "E11" = "  -3.63%  "
    "D12" = "3.939.70"
    "E12" = "  -0.68%  "
    "D13" = "0.125"
    "E13" = "  +1.86%  "
    "D14" = "27.44"
    "E14" = "  -2.81%  "
    "D15" = "3.372.52"
    "E15" = "  -0.57%  "
    "E16" = "  -2.84%  "
    "D17" = "60.923.06"
    "E17" = "  +0.08%  "
    "D18" = "6.03"
    "E18" = "  -3.52%  "
    "D19" = "13.48"
    "E19" = "  -3.84%  "
    "D20" = "8.79"
    "E20" = "  -2.68%  "
    "D21" = "380.99"
    "E21" = "  -1.57%  "
    "D22" = "75.15"
    "E22" = "  +1.96%  "
    "D23" = "0.545"
    "E23" = "  -2.63%  "
    "E24" = "  -0.08%  "
    "D25" = "0.0000110"
    "E25" = "  -5.97%  "
    "E26" = "  +5.57%  "
    "E27" = "  -0.02%  "
    "D28" = "7.10"
    "E28" = "  -4.21%  "
    "D29" = "7.81"
    "E29" = "  -2.08%  "
    "D30" = "2.10"
    "E30" = "  -2.13%  "
    "D32" = "1.33"
    "E32" = "  -5.16%  "
    "D33" = "22.72"
    "E33" = "  -4.03%  "
    "D34" = "6.82"
    "E34" = "  -1.94%  "
    "D35" = "165.96"
    "E35" = "  -0.67%  "
    "D36" = "4.88"
    "E36" = "  -1.75%  "
    "D37" = "3.403.28"
    "E37" = "  -0.40%  "
    "D38" = "1.43"
    "E38" = "  -4.11%  "
    "D39" = "0.0754"
    "E39" = "  -2.76%  "
    "B40" = "Mantle"
    "C40" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D40" = "0.769"
    "E40" = "  -1.67%  "
    "B41" = "EnergySwap"
    "C41" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D41" = "25.02"
    "E41" = "  -10.44%  "
    "D42" = "4.31"
    "E42" = "  -2.73%  "
    "D43" = "1.61"
    "E43" = "  -4.22%  "
    "D44" = "1.09"
    "E44" = "  -3.10%  "
    "D45" = "2.438.11"
    "E45" = "  -4.55%  "
    "E46" = "  -0.03%  "
    "D47" = "6.56"
    "E47" = "  -4.07%  "
    "D48" = "22.09"
    "E48" = "  -5.35%  "
    "D49" = "0.0257"
    "E49" = "  -4.76%  "
    "D50" = "2.01"
    "E50" = "  -3.07%  "
    "E51" = "  -3.57%  "
}

# These cells hold price text that Excel would otherwise auto-parse as a
# number (and mangle, e.g. "0.0000110" -> 1.1E-05), so force Text format,
# assign, then restore the original cell style.
$forceTextCells = @("D5", "D6", "D9", "D11", "D13", "D14", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50")

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    if ($forceTextCells -contains $addr) {
        $origStyle = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $updates[$addr]
        $rng.Style = $origStyle
    } else {
        $rng.Value = $updates[$addr]
    }
}
